$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.541.12"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.584.45"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "3.050.25"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "63.392.38"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000155"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.23%  "
$ws.Range("D17").Value = "2.592.08"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "555.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").Value = "0.0₃0865"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.415"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.629"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "0.0₆0235"
$ws.Range("E51").Value = "  +17.43%  "
